# Zeitblatt auf Stand - add October entries for rows 27, 29-33 (Oktober sheet)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Oktober")

# --- set cell text in the same order the author typed the diary, so new
# shared-string entries land at the same indices as the target workbook ---

# Row 27: 1h, prototype parts & pinmapping research
$ws.Range("C27").Value = "Teile gesucht um ein prototyp für das Pinmapping zubauen, wie auch nach lesen, prüfen was nicht stimmt"
$ws.Range("B27").Value = 1

# Row 29: soldering / cable break
$ws.Range("C29").Value = "Löten der Platine, Kabelbruch im Kabel festgestellt somit alternative überlegt"

# Row 30: 4h, re-soldering without ribbon cable
$ws.Range("C30").Value = "Umlöte, da kein Flachbandkabel verwendet wurde, Test gemacht ohne erfolg bekomme beim Transceive immer 0 zurück "
$ws.Range("B30").Value = 4

# Row 33: 1h, second prototype soldered
$ws.Range("C33").Value = "zweiten Prototyp gelötet, die Prototyp fertig gestellt. Idee für befestigung des CC265o mit kleinen adapter Platinen als gummi laschen"
$ws.Range("B33").Value = 1

# Row 31: 4h, cc2650 i2c tests
$ws.Range("C31").Value = "Test mit cc2650 i2c zum laufen zu bringen, Recherche in tests ohne erfolg "
$ws.Range("B31").Value = 4

# Row 32: 2.5h, tiva c + cc2650 oscilloscope test
$ws.Range("C32").Value = "Test mit tiva c und cc2650 am osziloskop, CC2650 vermutlich hw fehler "
$ws.Range("B32").Value = 2.5

# Wrap the detail column for the long entries (row 29 stays unwrapped,
# matching the author's original formatting)
$ws.Range("C27").WrapText = $true
$ws.Range("C30").WrapText = $true
$ws.Range("C31").WrapText = $true
$ws.Range("C32").WrapText = $true
$ws.Range("C33").WrapText = $true

# Row 34 has no text yet, but already carries the wrapped-text style
$ws.Range("C34").WrapText = $true

# Rows that now wrap onto two lines get taller
$ws.Rows.Item(27).RowHeight = 31.5
$ws.Rows.Item(30).RowHeight = 31.5
$ws.Rows.Item(33).RowHeight = 31.5

# Cursor moved on to the next empty detail cell
$ws.Range("C32").Select()

$ws.Calculate()
